$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I24").Value = 4.555243902439025

$rf = 8.11048780487805
for ($r = 25; $r -le 38; $r++) {
    $ws.Cells.Item($r, 9).Value = $rf
}
